$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.409.32'
$ws.Range('E2').Value = '  +3.61%  '
$ws.Range('D3').Value = '2.403.80'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.66%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('D9').Value = '2.428.58'
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('E13').Value = '  +3.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.38'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.85%  '
$ws.Range('E15').Value = '  +5.14%  '
$ws.Range('D16').Value = '2.883.70'
$ws.Range('E16').Value = '  +1.66%  '
$ws.Range('D17').Value = '62.311.24'
$ws.Range('E17').Value = '  +3.70%  '
$ws.Range('D18').Value = '2.436.65'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.89'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.83'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('E23').Value = '  +11.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.997'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('B25').Value = 'BabyDogeCoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D25').Value = '0.0₆0628'
$ws.Range('E25').Value = '  +116.75%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.44'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '620.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +10.34%  '
$ws.Range('E28').Value = '  +10.77%  '
$ws.Range('E29').Value = '  +4.95%  '
$ws.Range('D30').Value = '0.0₃0973'
$ws.Range('E30').Value = '  +5.18%  '
$ws.Range('D31').Value = '2.557.29'
$ws.Range('E31').Value = '  +2.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.13'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').Value = '  +6.80%  '
$ws.Range('E34').Value = '  +3.46%  '
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.995'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.89'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('E43').Value = '  +13.50%  '
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Polygon'
$ws.Range('C46').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.835'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +118.85%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '14.69'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +25.49%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '143.78'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.58'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.41'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.596'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.48%  '
